# Fix wrong Excel formatting: "Tin 5" (with stray space) -> "Tin5"
# This also causes the now-unused "Tin 5" shared-string entry to be
# dropped from the workbook on save, which re-packs the shared string
# table (the other cells referencing strings that sorted after the
# removed one shift down by one index, with their displayed text
# unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = "Tin5"

# Update the last selected cell to match the author's final selection.
$ws.Range("F10").Select()
